$p = $ppt.ActivePresentation

# -------------------------------------------------------------------
# Edit 1: Slide 6 speaker notes - remove the leftover question text,
# leaving an empty paragraph (as in the target diff).
# -------------------------------------------------------------------
$slide6 = $p.Slides.Item(6)
$notes6 = $slide6.NotesPage
$notesBody = $notes6.Shapes.Item("Notes Placeholder 2")
$notesBody.TextFrame.TextRange.Text = ""

# -------------------------------------------------------------------
# Edit 2: Slide 8 ("Iteration 1") body - split the bullet
# "Ported code from Iteration 0 into project" into two runs:
#   "Ported relevant " + "code from Iteration 0 into project"
# -------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$content8 = $slide8.Shapes.Item("Content Placeholder 2")
$tr = $content8.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
$targetIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text.Contains("Ported code from Iteration 0 into project")) {
        $targetIndex = $i
    }
}

$targetPara = $tr.Paragraphs($targetIndex, 1)
$firstRun = $targetPara.Runs(1)
$firstRun.Text = "code from Iteration 0 into project"
$targetPara.InsertBefore("Ported relevant ")
